# Modified testcases and log file
# - Column B (EXECUTE) flips several rows from "No" to "YES"
# - Column H (PARALLEL) flips rows 2-14 to the new value "Yes"
# - The sheet's active selection grows from H2:H13 to H2:H14

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B: rows that change from "No" to "YES"
$rowsYes = 2,4,5,12,13
foreach ($r in $rowsYes) {
    $ws.Range("B$r").Value = "YES"
}

# Column H: rows 2-14 all become "Yes" (new shared string)
for ($r = 2; $r -le 14; $r++) {
    $ws.Range("H$r").Value = "Yes"
}

# Update the sheet's visible selection to match the widened range
$ws.Range("H2:H14").Select()
